# Update Name of Algo
# Apply updated RandomForest imputation values to columns A and B.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "A11" = -21.7779
    "B11" = 5.384700000000001
    "A12" = -21.39029999999999
    "A15" = -21.68779999999999
    "B23" = 8.8193
    "A27" = -21.9692
    "A28" = -21.95510000000001
    "B28" = 5.526400000000001
    "A31" = -21.71940000000001
    "A32" = -21.05619999999999
    "B32" = 6.043299999999997
    "B34" = 9.437700000000008
    "A36" = -21.28179999999999
    "B36" = 5.584199999999998
    "B37" = 8.972600000000005
    "A38" = -19.69119999999998
    "B42" = 10.2315
    "A46" = -21.96270000000001
    "B49" = 4.985499999999999
    "A54" = -21.6992
    "B54" = 4.684899999999998
    "A55" = -21.88180000000001
    "A56" = -22.04160000000001
    "A67" = -21.45629999999997
    "A69" = -21.73629999999997
    "A72" = -22.0473
    "A73" = -19.2719
    "B78" = 9.408499999999997
    "B80" = 9.468200000000001
    "A83" = -21.60829999999999
    "A86" = -21.99560000000001
    "A91" = -20.45969999999998
    "A93" = -21.5235
    "B97" = 6.189199999999997
    "A99" = -21.8138
    "B99" = 5.580399999999998
    "B100" = 4.7805
    "B101" = 4.4384
    "A104" = -21.55279999999999
    "A105" = -19.98219999999999
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
